# Edit script for montenegro_prva-crnogorska-liga_2023-2024.xlsx
# 1) Swap the match data (columns F:V) between several pairs of existing
#    rows (the "home/away" ordering for those fixtures was corrected).
# 2) Append 6 new match rows (90-95) at the bottom of the sheet and extend
#    the used range / dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Part 1: swap columns F..V between paired rows
# ---------------------------------------------------------------------
$colFirst = 6   # column F
$colLast  = 22  # column V

function Get-RowVals($row) {
    $vals = @{}
    for ($col = $colFirst; $col -le $colLast; $col++) {
        $vals[$col] = $ws.Cells.Item($row, $col).Value2
    }
    return $vals
}

function Set-RowVals($row, $vals) {
    for ($col = $colFirst; $col -le $colLast; $col++) {
        $ws.Cells.Item($row, $col).Value2 = $vals[$col]
    }
}

$swapPairs = @(
    ,@(18, 20)
    ,@(24, 25)
    ,@(41, 42)
    ,@(47, 48)
    ,@(72, 73)
    ,@(77, 78)
    ,@(82, 83)
)

foreach ($pair in $swapPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]
    $valsA = Get-RowVals $rowA
    $valsB = Get-RowVals $rowB
    Set-RowVals $rowA $valsB
    Set-RowVals $rowB $valsA
}

# ---------------------------------------------------------------------
# Part 2: append the 6 new rows (90-95)
# ---------------------------------------------------------------------
$newRows = @(
    @{ row=90; A=89; E=45264.625; F='Mornar Bar'; G=2; H='Jezero'; I=1; J=2.4; K='01/12/2023 03:12'; L=2.65; M='04/12/2023 14:58'; N=2.73; O='01/12/2023 03:12'; P=2.5; Q='04/12/2023 14:58'; R=3.09; S='01/12/2023 03:12'; T=3.39; U='04/12/2023 14:58'; V='https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/mornar-bar-jezero/fchsCFzH/' }
    @{ row=91; A=90; E=45269.54166666666; F='Rudar'; G=2; H='Arsenal Tivat'; I=0; J=2.57; K='08/12/2023 01:12'; L=2.7; M='09/12/2023 12:59'; N=2.81; O='08/12/2023 01:12'; P=2.84; Q='09/12/2023 12:59'; R=2.77; S='08/12/2023 01:12'; T=2.87; U='09/12/2023 12:59'; V='https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/rudar-arsenal-tivat/zVSFQEkA/' }
    @{ row=92; A=91; E=45269.625; F='Decic'; G=3; H='Jedinstvo'; I=1; J=1.33; K='09/12/2023 13:34'; L=1.33; M='09/12/2023 13:34'; N=4.62; O='09/12/2023 13:34'; P=4.62; Q='09/12/2023 13:34'; R=7.69; S='09/12/2023 13:34'; T=7.69; U='09/12/2023 13:34'; V='https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/decic-jedinstvo/dvUBRfz4/' }
    @{ row=93; A=92; E=45269.66666666666; F='Petrovac'; G=0; H='Mornar Bar'; I=0; J=2.2; K='08/12/2023 04:12'; L=2.51; M='09/12/2023 15:38'; N=2.9; O='08/12/2023 04:12'; P=2.81; Q='09/12/2023 15:38'; R=3.25; S='08/12/2023 04:12'; T=3.16; U='09/12/2023 15:38'; V='https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/petrovac-mornar-bar/OCG4TG5i/' }
    @{ row=94; A=93; E=45270.54166666666; F='Mladost DG'; G=1; H='Buducnost'; I=3; J=5.17; K='09/12/2023 01:12'; L=5.61; M='10/12/2023 12:54'; N=3.69; O='09/12/2023 01:12'; P=3.82; Q='10/12/2023 12:54'; R=1.56; S='09/12/2023 01:12'; T=1.58; U='10/12/2023 12:54'; V='https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/mladost-dg-buducnost/23F8SzLc/' }
    @{ row=95; A=94; E=45270.625; F='Jezero'; G=0; H='Sutjeska'; I=0; J=3.68; K='09/12/2023 03:13'; L=3.72; M='10/12/2023 14:55'; N=2.98; O='09/12/2023 03:13'; P=2.94; Q='10/12/2023 14:55'; R=2; S='09/12/2023 03:13'; T=2.16; U='10/12/2023 14:55'; V='https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/jezero-sutjeska/IgIKPY4G/' }
)

# Copy the formatting of an existing data row (row 2: bold/bordered index
# cell in column A, date-time numeric format in column E, plain cells
# elsewhere) down across the whole new block in one shot.
$ws.Range("A2:V2").Copy()
$ws.Range("A90:V95").PasteSpecial(-4122)

foreach ($nr in $newRows) {
    $r = $nr.row
    $ws.Cells.Item($r, 1).Value2  = $nr.A          # Indice
    $ws.Cells.Item($r, 2).Value2  = "montenegro"   # pais
    $ws.Cells.Item($r, 3).Value2  = "prva-crnogorska-liga"  # torneio
    $ws.Cells.Item($r, 4).Value2  = "2023-2024"    # temporada
    $ws.Cells.Item($r, 5).Value2  = $nr.E          # data_partida
    $ws.Cells.Item($r, 6).Value2  = $nr.F          # home
    $ws.Cells.Item($r, 7).Value2  = $nr.G          # home_ft_gols
    $ws.Cells.Item($r, 8).Value2  = $nr.H          # away
    $ws.Cells.Item($r, 9).Value2  = $nr.I          # away_ft_gols
    $ws.Cells.Item($r, 10).Value2 = $nr.J          # home_opening_odds
    $ws.Cells.Item($r, 11).Value2 = $nr.K          # home_opening_data_hora
    $ws.Cells.Item($r, 12).Value2 = $nr.L          # home_closing_odds
    $ws.Cells.Item($r, 13).Value2 = $nr.M          # home_closing_data_hora
    $ws.Cells.Item($r, 14).Value2 = $nr.N          # draw_opening_odds
    $ws.Cells.Item($r, 15).Value2 = $nr.O          # draw_opening_data_hora
    $ws.Cells.Item($r, 16).Value2 = $nr.P          # draw_closing_odds
    $ws.Cells.Item($r, 17).Value2 = $nr.Q          # draw_closing_data_hora
    $ws.Cells.Item($r, 18).Value2 = $nr.R          # away_opening_odds
    $ws.Cells.Item($r, 19).Value2 = $nr.S          # away_opening_data_hora
    $ws.Cells.Item($r, 20).Value2 = $nr.T          # away_closing_odds
    $ws.Cells.Item($r, 21).Value2 = $nr.U          # away_closing_data_hora
    $ws.Cells.Item($r, 22).Value2 = $nr.V          # url_partida
}

Write-Host "Edit complete: swapped $($swapPairs.Count) row pairs, added $($newRows.Count) new rows."
